$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns before the old "pair" column (old B), shifting
#    old B..K (pair..random) right to D..M. This creates room for the new
#    "sort" (B) and "part" (C) columns.
# ---------------------------------------------------------------------------
$ws.Range("B:C").Insert()

# Headers for the two newly inserted columns. "part" must be interned into
# the shared-string table before "sort" (index 11 then 12) to match target.
$ws.Range("C1").Value = "part"
$ws.Range("B1").Value = "sort"

# ---------------------------------------------------------------------------
# 2. Fill in "sort" (B) and "part" (C) for the original 6 rows (part 1).
# ---------------------------------------------------------------------------
for ($i = 2; $i -le 7; $i++) {
    $ws.Cells.Item($i, 2).Value = $i - 1
    $ws.Cells.Item($i, 3).Value = 1
}

# ---------------------------------------------------------------------------
# 3. Part 2 rows (8-13): trial, sort, part, pair values.
#    (D = pair, pulled from the values the old sheet used for rows 8-13)
# ---------------------------------------------------------------------------
$part2Trial = @{8=7; 9=9; 10=11; 11=8; 12=10; 13=12}
$part2Pair  = @{8=12; 9=23; 10=31; 11=45; 12=56; 13=64}
foreach ($r in 8..13) {
    $ws.Cells.Item($r, 1).Value = $part2Trial[$r]
    $ws.Cells.Item($r, 2).Value = $r - 1
    $ws.Cells.Item($r, 3).Value = 2
    $ws.Cells.Item($r, 4).Value = $part2Pair[$r]
}

# ---------------------------------------------------------------------------
# 4. Part 3 rows (14-19): new rows entirely.
# ---------------------------------------------------------------------------
$part3Trial = @{14=17; 15=15; 16=13; 17=18; 18=16; 19=14}
$part3Pair  = @{14=14; 15=25; 16=36; 17=15; 18=26; 19=34}
foreach ($r in 14..19) {
    $ws.Cells.Item($r, 1).Value = $part3Trial[$r]
    $ws.Cells.Item($r, 2).Value = $r - 1
    $ws.Cells.Item($r, 3).Value = 3
    $ws.Cells.Item($r, 4).Value = $part3Pair[$r]
}

# ---------------------------------------------------------------------------
# 5. leftright_counter (K) / longdelay (L) for the new part-2 and part-3 rows.
# ---------------------------------------------------------------------------
$kVals = @{8=0; 9=1; 10=0; 11=1; 12=0; 13=1; 14=0; 15=1; 16=0; 17=1; 18=0; 19=1}
$lVals = @{8=1; 9=1; 10=1; 11=0; 12=0; 13=0; 14=0.5; 15=0.5; 16=0.5; 17=0.5; 18=0.5; 19=0.5}
foreach ($r in 8..19) {
    $ws.Cells.Item($r, 11).Value = $kVals[$r]
    $ws.Cells.Item($r, 12).Value = $lVals[$r]
}

# ---------------------------------------------------------------------------
# 6. stim_1 / stim_2 (E/F) formulas for the derived rows (part 2 and 3),
#    referencing the stim_1/stim_2 values of the original 6 pairs. These
#    carry the new "yellow highlight" style applied by the author.
# ---------------------------------------------------------------------------
$ws.Range("E8").Formula = "=E2"
$ws.Range("F8").Formula = "=F3"
$ws.Range("E9").Formula = "=E3"
$ws.Range("F9").Formula = "=F4"
$ws.Range("E10").Formula = "=E4"
$ws.Range("F10").Formula = "=F2"
$ws.Range("E11").Formula = "=E5"
$ws.Range("F11").Formula = "=F6"
$ws.Range("E12").Formula = "=E6"
$ws.Range("F12").Formula = "=F7"
$ws.Range("E13").Formula = "=E7"
$ws.Range("F13").Formula = "=F5"

$ws.Range("E14").Formula = "=E2"
$ws.Range("F14").Formula = "=E5"
$ws.Range("E15").Formula = "=E3"
$ws.Range("F15").Formula = "=E6"
$ws.Range("E16").Formula = "=E4"
$ws.Range("F16").Formula = "=E7"
$ws.Range("E17").Formula = "=F2"
$ws.Range("F17").Formula = "=F6"
$ws.Range("E18").Formula = "=F3"
$ws.Range("F18").Formula = "=F7"
$ws.Range("E19").Formula = "=F4"
$ws.Range("F19").Formula = "=F5"

$ws.Range("E8:F19").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 7. stim_left / stim_right / choice_stim_left / choice_stim_right formulas,
#    filled across the full A2:M19 grid (shared formulas, like the original
#    sheet already used for "random").
# ---------------------------------------------------------------------------
$ws.Range("G2:G19").Formula = "=IF(K2=0,E2,F2)"
$ws.Range("H2:H19").Formula = "=IF(K2=0,F2,E2)"
$ws.Range("I2:I19").Formula = "=IF(K2=0,CONCATENATE(""images/choice_trial_"",E2,"".png""),CONCATENATE(""images/choice_trial_"",F2,"".png""))"
$ws.Range("J2:J19").Formula = "=IF(K2=0,CONCATENATE(""images/choice_trial_"",F2,"".png""),CONCATENATE(""images/choice_trial_"",E2,"".png""))"

# ---------------------------------------------------------------------------
# 8. random (M) — keep row 2 and row 8 as independent RAND() cells (the
#    anchors for each part), and fill M3:M19 (skipping M8) as one shared
#    volatile formula, mirroring how the original K3:K13 group worked.
# ---------------------------------------------------------------------------
$ws.Range("M2").Formula = "=RAND()"
$ws.Range("M3:M7").Formula = "=RAND()"
$ws.Range("M8").Formula = "=RAND()"
$ws.Range("M9:M19").Formula = "=RAND()"

# ---------------------------------------------------------------------------
# 9. Re-apply the sort (by "pair", column D) over the original 6 rows, same
#    as the workbook's stored sortState.
# ---------------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D2:D7"))
$ws.Sort.SetRange($ws.Range("A2:M7"))
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 10. Misc view/metadata touch-ups to mirror the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("D4").Select()

